$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.382.95'
$ws.Range("E2").Value = '  +1.39%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.778.40'
$ws.Range("E3").Value = '  +3.72%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.000'
$ws.Range("E4").Value = '  -0.13%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '313.82'
$ws.Range("E5").Value = '  +1.65%  '
$ws.Range("E6").Value = '  -0.06%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5355'
$ws.Range("E7").Value = '  +12.81%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3767'
$ws.Range("E8").Value = '  +8.63%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '42.91'
$ws.Range("E9").Value = '  +2.14%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07402'
$ws.Range("E10").Value = '  +2.32%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.095'
$ws.Range("E11").Value = '  +5.39%  '
$ws.Range("E12").Value = '  -0.13%  '
$ws.Range("E13").Value = '  +4.55%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.103'
$ws.Range("E14").Value = '  +4.73%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '1.782.23'
$ws.Range("E15").Value = '  +3.42%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.993'
$ws.Range("E16").Value = '  +2.35%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '89.60'
$ws.Range("E17").Value = '  +3.22%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001057'
$ws.Range("E18").Value = '  +2.10%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06438'
$ws.Range("E19").Value = '  +1.09%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.9998'
$ws.Range("E20").Value = '  -0.11%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '16.80'
$ws.Range("E21").Value = '  +2.00%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.902'
$ws.Range("E22").Value = '  +5.26%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '27.415.69'
$ws.Range("E23").Value = '  +1.27%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.22'
$ws.Range("E24").Value = '  +4.61%  '
$ws.Range("E25").Value = '  -0.33%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '155.43'
$ws.Range("E26").Value = '  +3.04%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.18'
$ws.Range("E27").Value = '  +1.21%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.359'
$ws.Range("E28").Value = '  +14.25%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.985.09'
$ws.Range("E29").Value = '  +3.50%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '121.22'
$ws.Range("E30").Value = '  +0.74%  '
$ws.Range("E31").Value = '  +5.66%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.1026'
$ws.Range("E32").Value = '  +12.51%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.583'
$ws.Range("E33").Value = '  +5.35%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.622'
$ws.Range("E34").Value = '  +0.63%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.02259'
$ws.Range("E35").Value = '  +4.01%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.05977'
$ws.Range("E36").Value = '  +2.29%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.926'
$ws.Range("E37").Value = '  +4.79%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2060'
$ws.Range("E38").Value = '  +3.23%  '
$ws.Range("E39").Value = '  +3.08%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '8.273'
$ws.Range("E40").Value = '  +10.83%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.6113'
$ws.Range("E41").Value = '  +2.64%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.425'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.136'
$ws.Range("E43").Value = '  +4.97%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '13.23'
$ws.Range("E44").Value = '  +3.27%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.5782'
$ws.Range("E45").Value = '  +4.02%  '
$ws.Range("E46").Value = '  +1.46%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '121.22'
$ws.Range("E47").Value = '  +2.02%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.891'
$ws.Range("E48").Value = '  +3.78%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.120'
$ws.Range("E49").Value = '  +0.64%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06726'
$ws.Range("E50").Value = '  +1.47%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '70.79'
